$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" sheet as the first sheet in the book.
# ------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $info.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$idCell = $info.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "5736"
$info.Cells.Item(2, 2).Value = "Zak Crawley"
$info.Cells.Item(2, 3).Value = "Right Handed"
$info.Cells.Item(2, 4).Value = "Right Arm Off Break"

$info.Range("A1").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Update the "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE,
#    replacing the full scorecard URLs with just the match codes.
# ------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$lastRow = $batting.Cells.Item(1, 1).End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $link = [string]$cell.Value()
    if ($link -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}
